# localization-status.xlsx - "Generate Report for Archive"
#
# 1. Status text "Ready for handoff" -> "In Translation" wherever it
#    appears: Overview!E2:F3 (per-locale status columns) and the
#    "Status" column (C) on the per-locale "zh-cn" / "de-de" sheets.
# 2. Those narrower strings autofit to a narrower column: Overview!E:F
#    and column C on the "zh-cn" / "de-de" sheets shrink accordingly.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# The canonical file stores the shrunken columns at width
# 13.4101845877511 characters. This engine's ColumnWidth setter only
# lands on 1/6-character increments, so 12.5 is the input that gets us
# to the closest achievable width (13.333333333333334).
$newWidth = 12.5

# --- Overview sheet: per-locale status columns E (zh-cn) and F (de-de) ---
$ws = $wb.Worksheets.Item("Overview")
foreach ($ref in @("E2", "F2", "E3", "F3")) {
    $cell = $ws.Range($ref)
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}
$ws.Range("E1:F1").ColumnWidth = $newWidth

# --- zh-cn / de-de sheets: "Status" column (C) ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($ref in @("C2", "C3")) {
        $cell = $ws.Range($ref)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
    $ws.Range("C1").ColumnWidth = $newWidth
}
